$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '74.952.07'
$ws.Range("E2").Value = '  +1.49%  '

# Row 3
$ws.Range("D3").Value = '2.825.55'
$ws.Range("E3").Value = '  +7.83%  '

# Row 4
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '188.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.09%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '595.15'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.47%  '

# Row 7
$ws.Range("E7").Value = '  -0.01%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.551'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.35%  '

# Row 9
$ws.Range("E9").Value = '  -3.25%  '

# Row 10
$ws.Range("D10").Value = '2.824.06'
$ws.Range("E10").Value = '  +7.76%  '

# Row 11
$ws.Range("E11").Value = '  -1.04%  '

# Row 12
$ws.Range("E12").Value = '  +3.56%  '

# Row 13
$ws.Range("E13").Value = '  +2.31%  '

# Row 14
$ws.Range("D14").Value = '3.345.40'
$ws.Range("E14").Value = '  +8.74%  '

# Row 15
$ws.Range("D15").Value = '74.878.10'
$ws.Range("E15").Value = '  +2.01%  '

# Row 16
$ws.Range("E16").Value = '  -0.32%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.95'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.38%  '

# Row 18
$ws.Range("D18").Value = '2.822.67'
$ws.Range("E18").Value = '  +7.46%  '

# Row 19
$ws.Range("E19").Value = '  -1.83%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.32'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.12%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '376.38'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.08%  '

# Row 22
$ws.Range("E22").Value = '  -0.41%  '

# Row 23
$ws.Range("E23").Value = '  -0.21%  '

# Row 24
$ws.Range("E24").Value = '  -0.12%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '70.81'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.67%  '

# Row 27
$ws.Range("E27").Value = '  +0.74%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.67'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.82%  '

# Row 29
$ws.Range("E29").Value = '  +11.54%  '

# Row 30
$ws.Range("E30").Value = '  -0.87%  '

# Row 31
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.39'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.41%  '

# Row 32
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '511.36'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.75%  '

# Row 33
$ws.Range("E33").Value = '  +1.82%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.79'
$ws.Range("D34").Style = "Normal"

# Row 35
$ws.Range("E35").Value = '  -0.07%  '

# Row 36
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '20.02'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.41%  '

# Row 37
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '161.96'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.37%  '

# Row 38
$ws.Range("E38").Value = '  -0.64%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.41'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.72%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '185.16'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +15.63%  '

# Row 41
$ws.Range("E41").Value = '  +0.04%  '

# Row 42
$ws.Range("E42").Value = '  +4.78%  '

# Row 43
$ws.Range("E43").Value = '  +2.89%  '

# Row 44
$ws.Range("E44").Value = '  +0.29%  '

# Row 45
$ws.Range("E45").Value = '  +3.22%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.97'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.56%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.34'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.15%  '

# Row 48
$ws.Range("E48").Value = '  -0.93%  '

# Row 49
$ws.Range("E49").Value = '  +8.85%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.71'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.42%  '

# Row 51
$ws.Range("E51").Value = '  +8.19%  '
